$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.974.42'
$ws.Range("E2").Value = '  -0.46%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.562.37'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.63'
$ws.Range("E5").Value = '  -0.09%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.490'
$ws.Range("E6").Value = '  -0.11%  '

$ws.Range("E7").Value = '  +0.17%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.10'
$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("E9").Value = '  -0.25%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0599'
$ws.Range("E10").Value = '  +1.97%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0856'
$ws.Range("E11").Value = '  -0.55%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.787.07'
$ws.Range("E12").Value = '  -0.05%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.564.01'
$ws.Range("E13").Value = '  -0.14%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.75'
$ws.Range("E14").Value = '  -0.16%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.519'
$ws.Range("E15").Value = '  -0.07%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '62.00'
$ws.Range("E16").Value = '  +0.06%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.982.74'
$ws.Range("E17").Value = '  -0.50%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0704'
$ws.Range("E18").Value = '  +1.10%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '216.28'
$ws.Range("E19").Value = '  -1.47%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.34'
$ws.Range("E20").Value = '  -0.20%  '

$ws.Range("E21").Value = '  +0.15%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.10'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.22'
$ws.Range("E23").Value = '  -0.65%  '

$ws.Range("E24").Value = '  -0.56%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.97'
$ws.Range("E25").Value = '  -1.24%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.62'

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.07'
$ws.Range("E27").Value = '  +0.56%  '

$ws.Range("E28").Value = '  +1.12%  '

$ws.Range("E29").Value = '  +0.09%  '

$ws.Range("E31").Value = '  +0.97%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.23'
$ws.Range("E32").Value = '  -0.19%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.12'
$ws.Range("E33").Value = '  +1.02%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.422.74'
$ws.Range("E34").Value = '  -1.64%  '

$ws.Range("E35").Value = '  +2.94%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.06'
$ws.Range("E36").Value = '  +10.22%  '

$ws.Range("E37").Value = '  +1.90%  '

$ws.Range("E38").Value = '  -0.43%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.535'
$ws.Range("E39").Value = '  +2.44%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.82'
$ws.Range("E40").Value = '  +1.43%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.807'

$ws.Range("E42").Value = '  +0.18%  '

$ws.Range("E43").Value = '  +2.00%  '

$ws.Range("E44").Value = '  +1.61%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.68'
$ws.Range("E45").Value = '  +0.47%  '

$ws.Range("E46").Value = '  -1.85%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.699.59'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.36'
$ws.Range("E48").Value = '  +0.56%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0519'
$ws.Range("E49").Value = '  -0.89%  '

$ws.Range("E50").Value = '  -0.69%  '

$ws.Range("E51").Value = '  +0.12%  '
